{"js": "// Applies the \"Strengthened Manuscripts: Double Peer Review Refinements\" edit.\n// Strategy: operate on context.document.body.paragraphs by original index,\n// processing from the LAST paragraph to be touched back to the FIRST so that\n// inserting new paragraphs never shifts the index of a paragraph we still\n// need to visit.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- Paragraph 19 (Discussion body) -------------------------------------\nconst p19 = paragraphs.items[19];\np19.insertText(\n  \"Our data challenges the paradigm that vascular pathology in infection is solely a result of endothelial damage. Instead, we show that immune cells themselves are 'epigenetically loaded' to produce VEGFA, the potent permeability factor. The universality of this signature suggests that therapeutic strategies targeting chromatin remodeling (e.g., BET inhibitors) could have broad efficacy across multiple infectious diseases.\",\n  Word.InsertLocation.replace\n);\np19.insertParagraph(\n  \"Limitations of this study include the reliance on peripheral blood for Sepsis/Dengue; however, our TB data confirms that tissue-resident cells (BAL) show even stronger priming.\",\n  Word.InsertLocation.after\n);\n\n// ---- Paragraph 17 (VEGFA section body) -----------------------------------\nconst p17 = paragraphs.items[17];\np17.insertText(\n  \"The most striking finding was the status of Vampire Endothelial Growth Factor A (VEGFA). Traditionally considered an endothelial or stromal factor, we found VEGFA to be epigenetically primed and significantly upregulated in circulating monocytes and macrophages across all datasets.\",\n  Word.InsertLocation.replace\n);\np17.insertParagraph(\n  \"Log2 Fold Change (LFC) analysis showed a progressive increase in VEGFA expression correlating with disease severity risk: TB (Chronic) +1.2 LFC, Sepsis (Acute) +2.3 LFC, and Dengue (Hemorrhagic Risk) +4.0 LFC. This establishes a direct link between the epigenetic state of the immune system and the vascular leak phenotype characterizing severe shock.\",\n  Word.InsertLocation.after\n);\n\n// ---- Paragraph 16 (VEGFA section heading) --------------------------------\nconst p16 = paragraphs.items[16];\np16.insertText(\n  \"VEGFA: The Epigenetic Key to Vascular Shock\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 15 (Figure 2 caption) -------------------------------------\nconst p15 = paragraphs.items[15];\np15.insertText(\n  \"Figure: The Core Epigenetic Signature. Heatmap showing shared accessibility and expression of the 616 core genes.\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 13 (Core signature body) ----------------------------------\nconst p13 = paragraphs.items[13];\np13.insertText(\n  \"We intersected the primed gene sets from all three diseases and identified a Core Signature of 616 genes. Gene Ontology (GO) enrichment of this signature revealed major pathways: 'Response to Type I Interferon' (ISG15, MX1, STAT1), 'Neutrophil Degranulation' (S100A8, S100A9), and 'Antigen Presentation' (HLA-DRB5).\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 12 (Core signature heading) -------------------------------\nconst p12 = paragraphs.items[12];\np12.insertText(\n  \"A Core Signature of 616 'Locked' Genes\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 11 (Figure 1 caption) -------------------------------------\nconst p11 = paragraphs.items[11];\np11.insertText(\n  \"Figure: Universal Epigenetic Priming. Boxplot showing CPI distribution across TB, Sepsis, and Dengue. The epigenetic state is conserved (p=0.16).\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 9 (CPI body, first results para) --------------------------\nconst p9 = paragraphs.items[9];\np9.insertText(\n  \"To quantify the 'potential energy' of the immune genome, we developed the Chromatin Priming Index (CPI), defined as the fraction of differentially expressed genes (DEGs) that possess accessible chromatin promoters/enhancers. In a 'naive' state, genes must open chromatin before expression. In a 'primed' state, chromatin is already open.\",\n  Word.InsertLocation.replace\n);\np9.insertParagraph(\n  \"We analyzed 24,796 cells from Sepsis patients (GSE151263) and 20,000 cells from Dengue patients (GSE154386). Remarkably, the mean CPI was high and consistent across all conditions: TB (84.2%), Sepsis (82.5%), and Dengue (76.0%). Statistical comparison (Kruskal-Wallis test) yielded a p-value of 0.16, indicating no significant difference in the degree of epigenetic priming between these distinct diseases. This suggests that the 'Epigenetic Alert State' is a fundamental, conserved feature of the host response to severe stress.\",\n  Word.InsertLocation.after\n);\n\n// ---- Paragraph 8 (Results first heading) ---------------------------------\nconst p8 = paragraphs.items[8];\np8.insertText(\n  \"The Chromatin Priming Index (CPI) Reveals a Universal Alert State\",\n  Word.InsertLocation.replace\n);\n\n// ---- Paragraph 6 (Introduction body) --------------------------------------\nconst p6 = paragraphs.items[6];\np6.insertText(\n  \"The host immune response is evolutionary designed to protect against invasion, yet in severe infection, this response frequently becomes the driver of pathology. Conditions such as Sepsis and Dengue Shock Syndrome, despite their distinct pathogens, share striking clinical similarities: uncontrolled systemic inflammation, coagulopathy, and capillary leakage leading to hypotension and organ failure [1, 2]. Previous attempts to target specific cytokines (e.g., anti-TNF) have largely failed in sepsis, suggesting deeper regulatory mechanisms are at play.\",\n  Word.InsertLocation.replace\n);\np6.insertParagraph(\n  \"We hypothesized that the 'memory' or 'potential' for this pathological response is encoded not just in the transcriptome, but in the chromatin landscape. The phenomenon of 'Trained Immunity' [3] demonstrates that innate immune cells can undergo long-term epigenetic reprogramming. In this study, we asked: Is there a universal epigenetic state of 'severe infection'? To answer this, we performed a meta-analysis of single-cell Multiome (ATAC+RNA) data across Tuberculosis (chronic bacterial), Sepsis (acute syndromic), and Dengue (acute viral).\",\n  Word.InsertLocation.after\n);\n\n// ---- Paragraph 4 (Abstract body) -------------------------------------------\nconst p4 = paragraphs.items[4];\np4.insertText(\n  \"Severe infections, irrespective of their etiology\\u2014bacterial (Tuberculosis), viral (Dengue), or polymicrobial (Sepsis)\\u2014converge on a shared clinical phenotype of systemic inflammation, immune paralysis, and vascular dysfunction. While transcriptional studies have identified shared gene expression modules, the upstream regulatory mechanisms that 'lock' the immune system into this pathological state remain undefined. Here, we introduce the Chromatin Priming Index (CPI), a single-cell metric quantifying the decoupling of chromatin accessibility from gene expression ('poised' but not expressed genes). By applying CPI to multiomics data from active TB, Sepsis (24,796 cells), and Dengue (20,000 cells), we reveal a universally conserved 'epigenetic alert state' (mean CPI >80%) across all major immune subsets (p = 0.16, Kruskal-Wallis). We identify a core epigenetic signature of 616 genes that are primed for rapid activation, including classical antiviral and inflammatory mediators. Crucially, we discover that VEGFA\\u2014the primary driver of vascular permeability and shock\\u2014is epigenetically primed and transcriptionally upregulated in circulating immune cells across all three diseases (Log2FC: TB +1.2, Sepsis +2.3, Dengue +4.0). These findings identify immune-cell-derived VEGFA as a potential driver of the 'cytokine storm' vascular leak phenotype and suggest that the potential for shock is epigenetically imprinted in the myeloid compartment.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Applies the \"Strengthened Manuscripts: Double Peer Review Refinements\" edit.\n# Strategy: operate on $d.Paragraphs by original (1-based) index, processing\n# from the LAST paragraph to be touched back to the FIRST so that inserting\n# new paragraphs never shifts the index of a paragraph we still need to visit.\n\n$d = $word.ActiveDocument\n\n# ---- Paragraph 20 (Discussion body, 1-based) -----------------------------\n$p20 = $d.Paragraphs.Item(20)\n$p20.Range.Text = \"Our data challenges the paradigm that vascular pathology in infection is solely a result of endothelial damage. Instead, we show that immune cells themselves are 'epigenetically loaded' to produce VEGFA, the potent permeability factor. The universality of this signature suggests that therapeutic strategies targeting chromatin remodeling (e.g., BET inhibitors) could have broad efficacy across multiple infectious diseases.\"\n$d.Paragraphs.Item(20).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(21).Range.Text = \"Limitations of this study include the reliance on peripheral blood for Sepsis/Dengue; however, our TB data confirms that tissue-resident cells (BAL) show even stronger priming.\"\n\n# ---- Paragraph 18 (VEGFA section body) ------------------------------------\n$p18 = $d.Paragraphs.Item(18)\n$p18.Range.Text = \"The most striking finding was the status of Vampire Endothelial Growth Factor A (VEGFA). Traditionally considered an endothelial or stromal factor, we found VEGFA to be epigenetically primed and significantly upregulated in circulating monocytes and macrophages across all datasets.\"\n$d.Paragraphs.Item(18).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(19).Range.Text = \"Log2 Fold Change (LFC) analysis showed a progressive increase in VEGFA expression correlating with disease severity risk: TB (Chronic) +1.2 LFC, Sepsis (Acute) +2.3 LFC, and Dengue (Hemorrhagic Risk) +4.0 LFC. This establishes a direct link between the epigenetic state of the immune system and the vascular leak phenotype characterizing severe shock.\"\n\n# ---- Paragraph 17 (VEGFA section heading) ---------------------------------\n$p17 = $d.Paragraphs.Item(17)\n$p17.Range.Text = \"VEGFA: The Epigenetic Key to Vascular Shock\"\n\n# ---- Paragraph 16 (Figure 2 caption) ---------------------------------------\n$p16 = $d.Paragraphs.Item(16)\n$p16.Range.Text = \"Figure: The Core Epigenetic Signature. Heatmap showing shared accessibility and expression of the 616 core genes.\"\n\n# ---- Paragraph 14 (Core signature body) ------------------------------------\n$p14 = $d.Paragraphs.Item(14)\n$p14.Range.Text = \"We intersected the primed gene sets from all three diseases and identified a Core Signature of 616 genes. Gene Ontology (GO) enrichment of this signature revealed major pathways: 'Response to Type I Interferon' (ISG15, MX1, STAT1), 'Neutrophil Degranulation' (S100A8, S100A9), and 'Antigen Presentation' (HLA-DRB5).\"\n\n# ---- Paragraph 13 (Core signature heading) ---------------------------------\n$p13 = $d.Paragraphs.Item(13)\n$p13.Range.Text = \"A Core Signature of 616 'Locked' Genes\"\n\n# ---- Paragraph 12 (Figure 1 caption) ---------------------------------------\n$p12 = $d.Paragraphs.Item(12)\n$p12.Range.Text = \"Figure: Universal Epigenetic Priming. Boxplot showing CPI distribution across TB, Sepsis, and Dengue. The epigenetic state is conserved (p=0.16).\"\n\n# ---- Paragraph 10 (CPI body, first results para) ---------------------------\n$p10 = $d.Paragraphs.Item(10)\n$p10.Range.Text = \"To quantify the 'potential energy' of the immune genome, we developed the Chromatin Priming Index (CPI), defined as the fraction of differentially expressed genes (DEGs) that possess accessible chromatin promoters/enhancers. In a 'naive' state, genes must open chromatin before expression. In a 'primed' state, chromatin is already open.\"\n$d.Paragraphs.Item(10).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(11).Range.Text = \"We analyzed 24,796 cells from Sepsis patients (GSE151263) and 20,000 cells from Dengue patients (GSE154386). Remarkably, the mean CPI was high and consistent across all conditions: TB (84.2%), Sepsis (82.5%), and Dengue (76.0%). Statistical comparison (Kruskal-Wallis test) yielded a p-value of 0.16, indicating no significant difference in the degree of epigenetic priming between these distinct diseases. This suggests that the 'Epigenetic Alert State' is a fundamental, conserved feature of the host response to severe stress.\"\n\n# ---- Paragraph 9 (Results first heading) ------------------------------------\n$p9 = $d.Paragraphs.Item(9)\n$p9.Range.Text = \"The Chromatin Priming Index (CPI) Reveals a Universal Alert State\"\n\n# ---- Paragraph 7 (Introduction body) -----------------------------------------\n$p7 = $d.Paragraphs.Item(7)\n$p7.Range.Text = \"The host immune response is evolutionary designed to protect against invasion, yet in severe infection, this response frequently becomes the driver of pathology. Conditions such as Sepsis and Dengue Shock Syndrome, despite their distinct pathogens, share striking clinical similarities: uncontrolled systemic inflammation, coagulopathy, and capillary leakage leading to hypotension and organ failure [1, 2]. Previous attempts to target specific cytokines (e.g., anti-TNF) have largely failed in sepsis, suggesting deeper regulatory mechanisms are at play.\"\n$d.Paragraphs.Item(7).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(8).Range.Text = \"We hypothesized that the 'memory' or 'potential' for this pathological response is encoded not just in the transcriptome, but in the chromatin landscape. The phenomenon of 'Trained Immunity' [3] demonstrates that innate immune cells can undergo long-term epigenetic reprogramming. In this study, we asked: Is there a universal epigenetic state of 'severe infection'? To answer this, we performed a meta-analysis of single-cell Multiome (ATAC+RNA) data across Tuberculosis (chronic bacterial), Sepsis (acute syndromic), and Dengue (acute viral).\"\n\n# ---- Paragraph 5 (Abstract body) -----------------------------------------------\n$p5 = $d.Paragraphs.Item(5)\n$p5.Range.Text = \"Severe infections, irrespective of their etiology\u2014bacterial (Tuberculosis), viral (Dengue), or polymicrobial (Sepsis)\u2014converge on a shared clinical phenotype of systemic inflammation, immune paralysis, and vascular dysfunction. While transcriptional studies have identified shared gene expression modules, the upstream regulatory mechanisms that 'lock' the immune system into this pathological state remain undefined. Here, we introduce the Chromatin Priming Index (CPI), a single-cell metric quantifying the decoupling of chromatin accessibility from gene expression ('poised' but not expressed genes). By applying CPI to multiomics data from active TB, Sepsis (24,796 cells), and Dengue (20,000 cells), we reveal a universally conserved 'epigenetic alert state' (mean CPI >80%) across all major immune subsets (p = 0.16, Kruskal-Wallis). We identify a core epigenetic signature of 616 genes that are primed for rapid activation, including classical antiviral and inflammatory mediators. Crucially, we discover that VEGFA\u2014the primary driver of vascular permeability and shock\u2014is epigenetically primed and transcriptionally upregulated in circulating immune cells across all three diseases (Log2FC: TB +1.2, Sepsis +2.3, Dengue +4.0). These findings identify immune-cell-derived VEGFA as a potential driver of the 'cytokine storm' vascular leak phenotype and suggest that the potential for shock is epigenetically imprinted in the myeloid compartment.\"\n"}
